$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"

$ws.Range("C3").Value = "-"
$ws.Range("F3").Value = "-"

$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "['MCT-3A-Automação Industrial', -, 'MCT-2A-Acionamentos Elétricos', -]"
$ws.Range("F4").Value = "-"

$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "['MCT-3A-Automação Industrial', -, 'MCT-2A-Acionamentos Elétricos', -]"
$ws.Range("F6").Value = "-"

$ws.Range("C7").Value = "[-, 'MCT-3A-Lab. Máquinas Elétricas', -]"
$ws.Range("D7").Value = "['MCT-3A-Automação Industrial', -, 'MCT-2A-Acionamentos Elétricos', -]"
$ws.Range("F7").Value = "-"

$ws.Range("C8").Value = "[-, 'MCT-3A-Lab. Máquinas Elétricas', -]"
$ws.Range("D8").Value = "['MCT-3A-Automação Industrial', -, 'MCT-2A-Acionamentos Elétricos', -]"
